# Commit: "use lower case except for the first letter in title"
#
# The paper title that previously read:
#   "Nutritionally induced obesity Is attenuated in transgenic mice
#    overexpressing plasminogen activator inhibitor-1"
# had a spurious capital "Is". Fix the casing to "is" (lower case, as for
# every other word except the first) on every sheet that references this
# title, and leave the sheet's selection / active-cell state the way the
# author left it in the saved workbook.

$wb = $excel.ActiveWorkbook

$oldTitle = "Nutritionally induced obesity Is attenuated in transgenic mice overexpressing plasminogen activator inhibitor-1"
$newTitle = "Nutritionally induced obesity is attenuated in transgenic mice overexpressing plasminogen activator inhibitor-1"

# adipocyte_diameter: title lives in C7
$ws1 = $wb.Worksheets.Item("adipocyte_diameter")
$ws1.Cells.Item(7, 3).Value2 = $newTitle
$ws1.Range("C7").Select()

# adipose_vessel_size: title lives in C4
$ws2 = $wb.Worksheets.Item("adipose_vessel_size")
$ws2.Cells.Item(4, 3).Value2 = $newTitle
$ws2.Range("C4").Select()

# adipose_vessel_density: title lives in C4 - this is the sheet that ends up
# active/selected in the saved workbook, so select it last.
$ws3 = $wb.Worksheets.Item("adipose_vessel_density")
$ws3.Cells.Item(4, 3).Value2 = $newTitle
$ws3.Range("C4").Select()
$ws3.Activate()
